$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    # Row 36 is not part of this update (its F column value is a malformed
    # date, and it has not yet started its decrement cycle)
    if ($row -eq 36) { continue }

    $eCell = $ws.Cells.Item($row, 5)   # column E = "remaining"
    $dCell = $ws.Cells.Item($row, 4)   # column D = "total days"
    $fCell = $ws.Cells.Item($row, 6)   # column F = "start date"

    $eVal = $eCell.Value2
    $dVal = $dCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq 1) {
        # Cycle complete: reset remaining to the full total, and roll the
        # start date forward by the total number of days.
        $eCell.Value2 = $dVal
        $fCell.Value2 = $fVal + $dVal
    } else {
        # Normal day: one fewer day remaining.
        $eCell.Value2 = $eVal - 1
    }
}
